# Alert msg validation added and unwanted code removed
# (data-result workbook update: bump the DepVal test-id reference)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewPay")

# Row 2's DepVal (column D) now points at the newer test id.
$ws.Range("D2").Value = "EAPV21-0169"
